$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> @{ column letter = new text value }
# Price (D), Volume 1h (E) and Hora (G) columns refreshed with latest scrape.
$updates = @{
    2 = @{ "D" = "273.27"; "E" = "-1.64%"; "G" = "3" }
    3 = @{ "D" = "26.67"; "E" = "-2.28%"; "G" = "3" }
    4 = @{ "D" = "4.898"; "E" = "2.11%"; "G" = "3" }
    5 = @{ "D" = "0.06329"; "E" = "1.26%"; "G" = "3" }
    6 = @{ "D" = "6.900"; "E" = "0.92%"; "G" = "3" }
    7 = @{ "D" = "3.359"; "E" = "2.92%"; "G" = "3" }
    8 = @{ "D" = "1.263"; "E" = "34.03%"; "G" = "3" }
    9 = @{ "D" = "0.8771"; "E" = "0.18%"; "G" = "3" }
    10 = @{ "D" = "0.1457"; "E" = "-0.19%"; "G" = "3" }
    11 = @{ "D" = "0.04956"; "E" = "-1.67%"; "G" = "3" }
    12 = @{ "D" = "0.07340"; "E" = "0.97%"; "G" = "3" }
    13 = @{ "D" = "0.03126"; "E" = "-0.90%"; "G" = "3" }
    14 = @{ "D" = "0.09047"; "E" = "0.12%"; "G" = "3" }
    15 = @{ "D" = "0.001575"; "E" = "1.86%"; "G" = "3" }
    16 = @{ "D" = "0.0006329"; "E" = "0.76%"; "G" = "3" }
    17 = @{ "D" = "0.006037"; "E" = "5.06%"; "G" = "3" }
    18 = @{ "D" = "3.463"; "E" = "-0.22%"; "G" = "3" }
    19 = @{ "D" = "2.272"; "E" = "1.65%"; "G" = "3" }
    20 = @{ "G" = "3" }
    21 = @{ "E" = "1.38%"; "G" = "3" }
    22 = @{ "D" = "3.913"; "E" = "1.72%"; "G" = "3" }
    23 = @{ "D" = "0.04414"; "E" = "2.39%"; "G" = "3" }
    24 = @{ "D" = "0.001179"; "E" = "0.31%"; "G" = "3" }
    25 = @{ "D" = "0.003695"; "E" = "-13.32%"; "G" = "3" }
    26 = @{ "E" = "0.29%"; "G" = "3" }
    27 = @{ "E" = "5.61%"; "G" = "3" }
    28 = @{ "G" = "3" }
    29 = @{ "G" = "3" }
    30 = @{ "G" = "3" }
    31 = @{ "G" = "3" }
    32 = @{ "G" = "3" }
    33 = @{ "G" = "3" }
    34 = @{ "G" = "3" }
    35 = @{ "G" = "3" }
    36 = @{ "G" = "3" }
    37 = @{ "G" = "3" }
    38 = @{ "G" = "3" }
    39 = @{ "G" = "3" }
    40 = @{ "D" = "0.04040"; "E" = "0.48%"; "G" = "3" }
    41 = @{ "D" = "0.006637"; "E" = "3.79%"; "G" = "3" }
    42 = @{ "D" = "0.1166"; "E" = "1.43%"; "G" = "3" }
    43 = @{ "D" = "0.002125"; "E" = "-3.79%"; "G" = "3" }
    44 = @{ "D" = "0.01181"; "E" = "-13.50%"; "G" = "3" }
    45 = @{ "D" = "0.00005326"; "E" = "3.78%"; "G" = "3" }
    46 = @{ "E" = "17.11%"; "G" = "3" }
    47 = @{ "D" = "0.02004"; "E" = "-32.89%"; "G" = "3" }
    48 = @{ "G" = "3" }
    49 = @{ "G" = "3" }
    50 = @{ "G" = "3" }
    51 = @{ "G" = "3" }
}

$colIndex = @{ "D" = 4; "E" = 5; "G" = 7 }

foreach ($row in $updates.Keys) {
    foreach ($col in $updates[$row].Keys) {
        $cell = $ws.Cells.Item([int]$row, $colIndex[$col])
        $cell.NumberFormat = "@"
        $cell.Value = $updates[$row][$col]
    }
}
